$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the easting/northing coordinates in row 2
$ws.Range("Q2").Value = 402370
$ws.Range("R2").Value = 6710618

# Clear the time cells (Starttid / Sluttid) entirely for row 2
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
